$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E16:E22) lists period codes (YYMM) for
# rows 16-22. The old periods (2412..2506) are being replaced: oldest
# period (2412) is dropped from the top and the list now runs newest
# first (2506 down to 2412), i.e. the block of 7 rows is reversed.
$periods = @("2506", "2505", "2504", "2503", "2502", "2501", "2412")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# The "Valor Mora" amounts in column F follow the same row reversal:
# row 16 now carries the amount that used to be on row 22, and vice
# versa (the other rows keep the same 138000 value either way).
$ws.Range("F16").Value = 119600
$ws.Range("F22").Value = 138000
